# Apply updated cryptos list values (text), preserving original cell
# formatting by restoring default formats after forcing a Text number
# format during assignment (prevents Excel auto-converting numeric-
# looking strings like "353.31" into real numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "52.116.00"
Set-TextValue "E2" "  +0.87%  "

Set-TextValue "D3" "3.003.17"
Set-TextValue "E3" "  +2.75%  "

Set-TextValue "D5" "353.31"
Set-TextValue "E5" "  -0.28%  "

Set-TextValue "D6" "106.71"
Set-TextValue "E6" "  -2.99%  "

Set-TextValue "E7" "  -0.59%  "

Set-TextValue "E8" "  +0.20%  "

Set-TextValue "E9" "  -4.22%  "

Set-TextValue "D10" "38.09"
Set-TextValue "E10" "  -2.53%  "

Set-TextValue "E11" "  +2.04%  "

Set-TextValue "D12" "0.0855"
Set-TextValue "E12" "  -4.17%  "

Set-TextValue "D13" "19.01"
Set-TextValue "E13" "  -3.23%  "

Set-TextValue "D14" "3.471.56"
Set-TextValue "E14" "  +2.68%  "

Set-TextValue "D15" "7.61"
Set-TextValue "E15" "  -4.14%  "

Set-TextValue "D16" "2.986.10"
Set-TextValue "E16" "  +2.69%  "

Set-TextValue "E17" "  +3.37%  "

Set-TextValue "D18" "52.111.78"
Set-TextValue "E18" "  +0.73%  "

Set-TextValue "D19" "3.43"
Set-TextValue "E19" "  +5.09%  "

Set-TextValue "D20" "7.48"
Set-TextValue "E20" "  -1.76%  "

Set-TextValue "D21" "13.51"
Set-TextValue "E21" "  -4.88%  "

Set-TextValue "D22" "0.0₃0969"
Set-TextValue "E22" "  -1.21%  "

Set-TextValue "E23" "  -2.56%  "

Set-TextValue "D24" "264.17"
Set-TextValue "E24" "  -2.12%  "

Set-TextValue "E25" "  -3.50%  "

Set-TextValue "D26" "0.177"
Set-TextValue "E26" "  -3.17%  "

Set-TextValue "D27" "26.84"
Set-TextValue "E27" "  -1.33%  "

Set-TextValue "E29" "  -0.78%  "

Set-TextValue "E30" "  -1.59%  "

Set-TextValue "D31" "6.37"
Set-TextValue "E31" "  +3.49%  "

Set-TextValue "E32" "  -3.69%  "

Set-TextValue "E33" "  +15.48%  "

Set-TextValue "D34" "35.85"
Set-TextValue "E34" "  -5.18%  "

Set-TextValue "D35" "51.12"
Set-TextValue "E35" "  -2.25%  "

Set-TextValue "D36" "0.0437"
Set-TextValue "E36" "  -0.72%  "

Set-TextValue "E37" "  -0.06%  "

Set-TextValue "E38" "  +1.92%  "

Set-TextValue "D39" "2.84"
Set-TextValue "E39" "  +3.72%  "

Set-TextValue "E40" "  -2.59%  "

Set-TextValue "D41" "17.56"
Set-TextValue "E41" "  -3.68%  "

Set-TextValue "D42" "0.117"
Set-TextValue "E42" "  -1.04%  "

Set-TextValue "B43" "Monero"
Set-TextValue "C43" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D43" "124.33"
Set-TextValue "E43" "  +7.43%  "

Set-TextValue "B44" "EnergySwap"
Set-TextValue "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "22.62"
Set-TextValue "E44" "  -2.23%  "

Set-TextValue "E45" "  -1.85%  "

Set-TextValue "D46" "2.119.76"
Set-TextValue "E46" "  -0.86%  "

Set-TextValue "D47" "3.32"
Set-TextValue "E47" "  -3.87%  "

Set-TextValue "E48" "  -6.17%  "

Set-TextValue "D49" "3.296.24"
Set-TextValue "E49" "  +2.64%  "

Set-TextValue "D50" "0.243"
Set-TextValue "E50" "  -2.57%  "

Set-TextValue "D51" "0.0329"
Set-TextValue "E51" "  -1.03%  "
